$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates ---
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 10

# --- Zero out the per-line "Pricing" / H column values that were billed ---
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0

$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("H27").Value = 0

$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0

# --- Insert a new row above the old "TOTAL" row (row 35), pushing TOTAL to row 36 ---
$ws.Rows.Item(35).Insert()

# Copy formatting (fill/alignment/number format) from the row-33 style band so the
# new row matches the alternating shaded-row look (styles 12/13/14) instead of the
# plain band it inherited from row 34 on insert.
$ws.Range("A33:H33").Copy()
$ws.Range("A35:H35").PasteSpecial(-4122)

# --- Row 34 becomes the new "Point 09 / ANC-DSC-16-96-D1" entry, zeroed out ---
$ws.Range("A34").Value = "Point 09"
$ws.Range("B34").Value = "ANC-DSC-16-96-D1"
$ws.Range("D34").Value = "ANC,Disc,16in,96in,Db Eye 1in"
$ws.Range("F34").Value = 0
$ws.Range("H34").Value = 0

# --- New row 35 restores the original "Point 05 / ANC-EXP-8-72-S58" entry, zeroed out ---
$ws.Range("A35").Value = "Point 05"
$ws.Range("B35").Value = "ANC-EXP-8-72-S58"
$ws.Range("C35").Value = "Inst"
$ws.Range("D35").Value = "ANC,Expanding,8in,72in,Sg Eye 5/8in"
$ws.Range("E35").Value = "EA"
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = ""
$ws.Range("H35").Value = 0

# --- Row 36 (former row 35) is the TOTAL row, now zeroed ---
$ws.Range("H36").Value = 0

Write-Output "Edit complete"
